$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H68").Value = 20295
$ws.Range("J68").Value = 20295
$ws.Range("L68").Value = 20295
$ws.Range("N68").Value = -21793
$ws.Range("H70").Value = 2386.8
$ws.Range("I70").Value = 2560.2
$ws.Range("J70").Value = 2040
$ws.Range("K70").Value = 7680.599999999999
$ws.Range("L70").Value = 6120
$ws.Range("M70").Value = -7410.599999999999
$ws.Range("N70").Value = -6660
$ws.Range("H71").Value = 20295
$ws.Range("J71").Value = 20295
$ws.Range("L71").Value = 60885
$ws.Range("N71").Value = -68373
$ws.Range("H73").Value = 2386.8
$ws.Range("I73").Value = 2560.2
$ws.Range("J73").Value = 2040
$ws.Range("K73").Value = 7680.599999999999
$ws.Range("L73").Value = 6120
$ws.Range("M73").Value = -6744.599999999999
$ws.Range("N73").Value = -7992
$ws.Range("H121").Value = 1630
$ws.Range("I121").Value = 890
$ws.Range("J121").Value = 2000
$ws.Range("K121").Value = 2670
$ws.Range("L121").Value = 6000
$ws.Range("M121").Value = -923
$ws.Range("N121").Value = -9494
$ws.Range("H131").Value = 1396.7273
$ws.Range("J131").Value = 2452.5
$ws.Range("L131").Value = 7357.5
$ws.Range("N131").Value = -17437.5
$ws.Range("H132").Value = 1808.2759
$ws.Range("I132").Value = 1393.0869
$ws.Range("J132").Value = 3399.8333
$ws.Range("K132").Value = 4179.2607
$ws.Range("L132").Value = 10199.4999
$ws.Range("M132").Value = -1649.2607
$ws.Range("N132").Value = -15259.4999
$ws.Range("H138").Value = 7283.5107
$ws.Range("I138").Value = 1495.9
$ws.Range("J138").Value = 8847.7295
$ws.Range("K138").Value = 4487.700000000001
$ws.Range("L138").Value = 26543.1885
$ws.Range("M138").Value = 652.2999999999993
$ws.Range("N138").Value = -36823.1885
$ws.Range("H141").Value = 2243.9285
$ws.Range("I141").Value = 2150.8333
$ws.Range("J141").Value = 2802.5
$ws.Range("K141").Value = 6452.499899999999
$ws.Range("L141").Value = 8407.5
$ws.Range("M141").Value = -1272.499899999999
$ws.Range("N141").Value = -18767.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4974.385
$ws.Range("I32").Value = 3552.838
$ws.Range("K32").Value = 3552.838
$ws.Range("M32").Value = -3265.838
$ws.Range("H61").Value = 1924.85
$ws.Range("I61").Value = 2078.9333
$ws.Range("J61").Value = 1462.6
$ws.Range("K61").Value = 2078.9333
$ws.Range("L61").Value = 1462.6
$ws.Range("M61").Value = -1866.9333
$ws.Range("N61").Value = -1886.6
$ws.Range("H102").Value = 4632882
$ws.Range("I102").Value = 9261262
$ws.Range("J102").Value = 4502.75
$ws.Range("K102").Value = 9261262
$ws.Range("L102").Value = 4502.75
$ws.Range("M102").Value = -9259640
$ws.Range("N102").Value = -7746.75
$ws.Range("H122").Value = 1604261
$ws.Range("I122").Value = 1833126.9
$ws.Range("J122").Value = 2200
$ws.Range("K122").Value = 5499380.699999999
$ws.Range("L122").Value = 6600
$ws.Range("M122").Value = -5496930.699999999
$ws.Range("N122").Value = -11500
$ws.Range("H136").Value = 1924.85
$ws.Range("I136").Value = 2078.9333
$ws.Range("J136").Value = 1462.6
$ws.Range("K136").Value = 6236.7999
$ws.Range("L136").Value = 4387.799999999999
$ws.Range("M136").Value = -3686.7999
$ws.Range("N136").Value = -9487.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H94").Value = 1358.8043
$ws.Range("I94").Value = 1190.625
$ws.Range("J94").Value = 1743.2142
$ws.Range("K94").Value = 1190.625
$ws.Range("L94").Value = 1743.2142
$ws.Range("M94").Value = -739.625
$ws.Range("N94").Value = -2645.2142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 472.0909
$ws.Range("I107").Value = 436.66666
$ws.Range("J107").Value = 514.6
$ws.Range("K107").Value = 436.66666
$ws.Range("L107").Value = 514.6
$ws.Range("M107").Value = 1483.33334
$ws.Range("N107").Value = -4354.6
$ws.Range("H122").Value = 1010.3333
$ws.Range("I122").Value = 1010.3333
$ws.Range("K122").Value = 3030.9999
$ws.Range("M122").Value = -580.9998999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 200001890
$ws.Range("I32").Value = 1000000000
$ws.Range("J32").Value = 2350
$ws.Range("K32").Value = 3000000000
$ws.Range("L32").Value = 7050
$ws.Range("M32").Value = -2999999717
$ws.Range("N32").Value = -7616

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 10114
$ws.Range("J109").Value = 10114
$ws.Range("L109").Value = 10114
$ws.Range("N109").Value = -12194

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3650
$ws.Range("I7").Value = 2800
$ws.Range("J7").Value = 4500
$ws.Range("K7").Value = 2800
$ws.Range("L7").Value = 4500
$ws.Range("M7").Value = -2688
$ws.Range("N7").Value = -4724
$ws.Range("H68").Value = 100001940
$ws.Range("I68").Value = 1396.6666
$ws.Range("J68").Value = 250002750
$ws.Range("K68").Value = 1396.6666
$ws.Range("L68").Value = 250002750
$ws.Range("M68").Value = -647.6666
$ws.Range("N68").Value = -250004248
$ws.Range("H71").Value = 100001940
$ws.Range("I71").Value = 1396.6666
$ws.Range("J71").Value = 250002750
$ws.Range("K71").Value = 6983.333000000001
$ws.Range("L71").Value = 1250013750
$ws.Range("M71").Value = -3239.333000000001
$ws.Range("N71").Value = -1250021238
$ws.Range("H93").Value = 8963.27
$ws.Range("I93").Value = 12107.167
$ws.Range("J93").Value = 1889.5
$ws.Range("K93").Value = 12107.167
$ws.Range("L93").Value = 1889.5
$ws.Range("M93").Value = -10859.167
$ws.Range("N93").Value = -4385.5
$ws.Range("H126").Value = 3650
$ws.Range("I126").Value = 2800
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 8400
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -5930
$ws.Range("N126").Value = -18440
$ws.Range("H132").Value = 20637458
$ws.Range("I132").Value = 28891062
$ws.Range("J132").Value = 3443.1667
$ws.Range("K132").Value = 86673186
$ws.Range("L132").Value = 10329.5001
$ws.Range("M132").Value = -86670656
$ws.Range("N132").Value = -15389.5001
$ws.Range("H136").Value = 8914.294
$ws.Range("I136").Value = 3967.6365
$ws.Range("J136").Value = 17983.166
$ws.Range("K136").Value = 11902.9095
$ws.Range("L136").Value = 53949.49800000001
$ws.Range("M136").Value = -9352.9095
$ws.Range("N136").Value = -59049.49800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1827.3043
$ws.Range("I113").Value = 1339.0625
$ws.Range("J113").Value = 2943.2856
$ws.Range("K113").Value = 4017.1875
$ws.Range("L113").Value = 8829.856800000001
$ws.Range("M113").Value = -1847.1875
$ws.Range("N113").Value = -13169.8568
$ws.Range("H132").Value = 1897.6857
$ws.Range("I132").Value = 1210.85
$ws.Range("J132").Value = 2813.4666
$ws.Range("K132").Value = 3632.55
$ws.Range("L132").Value = 8440.399800000001
$ws.Range("M132").Value = -1102.55
$ws.Range("N132").Value = -13500.3998

